$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A3:K3")
$rng.NumberFormat = "@"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " November 02 2020"
$ws.Range("C3").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "Shahbaz Ahmed "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "100.00"
